$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present in the target sheet:
#   row 6: "free games helix jump " / "com.singleton.helix" (with styled, blank C6/D6)
#   row 7: duplicate "blockchain" / "block.chain.technology"
# Deleting both shifts rows 8-10 up to become rows 6-8.
$ws.Rows("6:7").Delete()

# Update the selected cell to match the post-edit selection state.
$ws.Range("A6").Select()
